$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# xlUp = -4162 : find the last used row in column A (Beteckning), walking up
# from the bottom of the sheet. This avoids relying on UsedRange, which can
# be thrown off by the sheet's leading blank row 0.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$newValue = 45208

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newValue
}
